# Adds two new "word" entries (bounty, distort) to the "word" sheet and
# two new phrase entries (be intended for, work toward) to the "phrase" sheet.

$wb = $excel.ActiveWorkbook

# --- "word" sheet: append bounty / distort after the last existing row (74) ---
$wordWs = $wb.Worksheets.Item("word")

$wordWs.Cells.Item(75, 1).Value = "bounty"
$wordWs.Cells.Item(75, 2).Value = "/ˈbaʊnti/"

$wordWs.Cells.Item(76, 1).Value = "distort"
$wordWs.Cells.Item(76, 2).Value = "/dɪˈstɔːrt/"

# --- "phrase" sheet: append "be intended for" and "work toward" after row 78 ---
$phraseWs = $wb.Worksheets.Item("phrase")

$phraseWs.Cells.Item(79, 1).Value = "be intended for"
$phraseWs.Cells.Item(79, 2).Value = "被设计给/被用来给/是为了…而准备的"
$phraseWs.Cells.Item(79, 3).Value = "表示某物的预期用途、目标用户或目的。"
$phraseWs.Cells.Item(79, 4).Value = "MySQL Server is intended for mission-critical, heavy-load production systems as well as for embedding into mass-deployed software."

$phraseWs.Cells.Item(80, 1).Value = "work toward"
$phraseWs.Cells.Item(80, 2).Value = "朝着某个方向努力、逐步实现"
$phraseWs.Cells.Item(80, 3).Value = "强调这是一个进行中的、有意识的目标，而不是已经完成的状态。"
$phraseWs.Cells.Item(80, 4).Value = "One of our main goals with the product is to continue to work toward compliance with the SQL standard, but without sacrificing speed or reliability."

# --- restore view state (scroll position / selection) to match the saved workbook ---
[void]$wordWs.Activate()
$excel.ActiveWindow.ScrollRow = 55
[void]$wordWs.Range("D80").Select()

[void]$phraseWs.Activate()
$excel.ActiveWindow.ScrollRow = 59
[void]$phraseWs.Range("D86").Select()

[void]$wordWs.Activate()
